$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new monthly data point (2024-11-01) was added to the top of the time
# series, pushing the existing rows down by one. Insert a fresh row 3
# (Excel copies formatting down from row 2, matching the target file).
$ws.Rows(3).Insert()

# New data point: 2024-11-01 (serial 45597), value 45900
$ws.Range("A3").Value2 = 45597
$ws.Range("B3").Value2 = 45900

# The "% change vs last year" column used to hold hard-coded literal
# values; it is now driven by a live formula comparing this row's value
# to the value 12 rows below (12 months prior).
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Formula = "=(B2/B14-1)*100"
$ws.Range("C2").Style = "Normal"

# Fill the formula down the rest of the column (C3:C33 is one shared
# formula group, C34:C65 a second one like the original authored file -
# the extra filled-but-unused rows are trimmed back out afterwards).
$ws.Range("C3:C33").Formula = "=(B3/B15-1)*100"
$ws.Range("C3:C33").Style = "Normal"

$ws.Range("C34:C65").Formula = "=(B34/B46-1)*100"
$ws.Range("C34:C65").Style = "Normal"

# Trim the filled-but-unused tail rows back out without disturbing any
# relative formula references into that range (ClearContents instead of
# a real row Delete, which would shift/#REF! the remaining formulas).
$ws.Rows("63:65").ClearContents()
$ws.Range("C62").ClearContents()
